$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# "Tip" column header becomes "Tip verificare"
$ws.Range("F1").Value = "Tip verificare"
# New column G header
$ws.Range("G1").Value = "Valori verificare"

# --- Update CodNFC (column D) numeric codes ---
$ws.Range("D2:D5").Value = 100003
$ws.Range("D6:D11").Value = 100004
$ws.Range("D12:D17").Value = 100005

# --- Row 16 count value change ---
$ws.Range("F16").Value = 3

# --- New "Valori verificare" column content ---
$ws.Range("G6").Value = "Functional, Nefuncțional"
$ws.Range("G12").Value = "Roșu, Maro, Verde"
$ws.Range("G16").Value = "Roșu, Maro, Verde"

# --- View state updates (matches final selection in the source workbook) ---
[void]$ws.Range("G16").Select()
